$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the partner-number cells (H/I) that reference the renumbered rows 42-45
$ws.Range("H6").Value = 43
$ws.Range("I6").Value = 43

$ws.Range("H11").Value = 42
$ws.Range("I11").Value = 42

$ws.Range("H20").Value = 44
$ws.Range("I20").Value = 44

$ws.Range("H38").Value = 41
$ws.Range("I38").Value = 41

# Renumber rows 42-45 (column A "Nr.") down by one
$ws.Range("A42").Value = 41
$ws.Range("A43").Value = 42
$ws.Range("A44").Value = 43
$ws.Range("A45").Value = 44

# Update the view: scroll so row 16 is the top-left visible row, and
# select H38 as the active cell
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("H38").Select()
